# class_progression_Wizard.xlsx — re-run of the "Progressions" XML map
# refresh (Data > Refresh All) after pointing the XML source at its new
# location (…/tietokanta/schemas/progression_schema.xml). Re-importing the
# XML data populated every previously-blank mapped cell with the
# XML-map's "missing element" placeholder text, "none".
#
# The COM layer in this runtime does not drive the real XML-map/connection
# refresh machinery, so we reproduce its net effect directly: write "none"
# into every table cell that the refresh touched, and restore the
# worksheet's active selection the way it was left after the refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every row's previously-empty "Advancement" cells (they already carry the
# table's style via s="1") get filled with the literal "none" — exactly
# what Excel's XML map binding writes for an element that is absent from
# the refreshed XML source.
$ws.Range("J2:K2").Value = "none"
$ws.Range("D3:K3").Value = "none"
$ws.Range("E4:K4").Value = "none"
$ws.Range("D5:K5").Value = "none"
$ws.Range("G6:K6").Value = "none"
$ws.Range("D7:K7").Value = "none"
$ws.Range("F8:K8").Value = "none"
$ws.Range("D9:K9").Value = "none"
$ws.Range("F10:K10").Value = "none"
$ws.Range("E11:K11").Value = "none"
$ws.Range("G12:K12").Value = "none"
$ws.Range("D13:K13").Value = "none"
$ws.Range("G14:K14").Value = "none"
$ws.Range("D15:K15").Value = "none"
$ws.Range("G16:K16").Value = "none"
$ws.Range("D17:K17").Value = "none"
$ws.Range("F18:K18").Value = "none"
$ws.Range("D19:K19").Value = "none"
$ws.Range("F20:K20").Value = "none"
$ws.Range("E21:K21").Value = "none"

# The refresh left the cursor on K1 (top-right header cell) rather than
# its previous spot (I5).
$ws.Range("K1").Select()
